# Apply scheduled market-price / profit recompute to the Leve profit tracker sheets.
# Updates currentAveragePrice(NQ/HQ) (H:K), derived Leve prices (L), and profit columns (M:N)
# for the rows whose source market data changed in this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1591.25
$ws.Range("I43").Value = 1626.909
$ws.Range("J43").Value = 1199
$ws.Range("K43").Value = 1626.909
$ws.Range("L43").Value = 1199
$ws.Range("M43").Value = -1557.909
$ws.Range("N43").Value = -1337
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2828
$ws.Range("H125").Value = 1349.1428
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1349.1428
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 12142.2852
$ws.Range("N125").Value = -17062.2852
$ws.Range("H127").Value = 1211.5
$ws.Range("I127").Value = 1255.8
$ws.Range("J127").Value = 990
$ws.Range("K127").Value = 3767.4
$ws.Range("L127").Value = 2970
$ws.Range("M127").Value = 1192.6
$ws.Range("N127").Value = -12890
$ws.Range("H135").Value = 1461.1428
$ws.Range("I135").Value = 1461.1428
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13150.2852
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -10615.2852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 14400
$ws.Range("I19").Value = 3800
$ws.Range("J19").Value = 25000
$ws.Range("K19").Value = 3800
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = -3571
$ws.Range("N19").Value = -25458
$ws.Range("H61").Value = 4391
$ws.Range("I61").Value = 4188.75
$ws.Range("J61").Value = 5200
$ws.Range("K61").Value = 4188.75
$ws.Range("L61").Value = 5200
$ws.Range("M61").Value = -3976.75
$ws.Range("N61").Value = -5624
$ws.Range("H110").Value = 5271
$ws.Range("I110").Value = 5617.1
$ws.Range("J110").Value = 4117.3335
$ws.Range("K110").Value = 5617.1
$ws.Range("L110").Value = 4117.3335
$ws.Range("M110").Value = -3572.1
$ws.Range("N110").Value = -8207.333500000001
$ws.Range("H131").Value = 78749.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 78749.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 78749.5
$ws.Range("N131").Value = -88829.5
$ws.Range("H132").Value = 2318.8865
$ws.Range("I132").Value = 2328.5527
$ws.Range("J132").Value = 2257.6667
$ws.Range("K132").Value = 6985.658100000001
$ws.Range("L132").Value = 6773.000100000001
$ws.Range("M132").Value = -4455.658100000001
$ws.Range("N132").Value = -11833.0001
$ws.Range("H136").Value = 4391
$ws.Range("I136").Value = 4188.75
$ws.Range("J136").Value = 5200
$ws.Range("K136").Value = 12566.25
$ws.Range("L136").Value = 15600
$ws.Range("M136").Value = -10016.25
$ws.Range("N136").Value = -20700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5323.647
$ws.Range("I105").Value = 5076.636
$ws.Range("J105").Value = 5776.5
$ws.Range("K105").Value = 5076.636
$ws.Range("L105").Value = 5776.5
$ws.Range("M105").Value = -3329.636
$ws.Range("N105").Value = -9270.5
$ws.Range("H107").Value = 2249.75
$ws.Range("I107").Value = 1999.6666
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 1999.6666
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -79.66660000000002
$ws.Range("N107").Value = -6840
$ws.Range("H134").Value = 4688.5557
$ws.Range("I134").Value = 5353.316
$ws.Range("J134").Value = 3109.75
$ws.Range("K134").Value = 16059.948
$ws.Range("L134").Value = 9329.25
$ws.Range("M134").Value = -13524.948
$ws.Range("N134").Value = -14399.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 607
$ws.Range("I5").Value = 618.5
$ws.Range("J5").Value = 599.3333
$ws.Range("K5").Value = 618.5
$ws.Range("L5").Value = 599.3333
$ws.Range("M5").Value = -506.5
$ws.Range("N5").Value = -823.3333
$ws.Range("H26").Value = 4071.2144
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 4071.2144
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 4071.2144
$ws.Range("N26").Value = -4645.2144
$ws.Range("H31").Value = 2974.353
$ws.Range("I31").Value = 1469.2307
$ws.Range("J31").Value = 7866
$ws.Range("K31").Value = 1469.2307
$ws.Range("L31").Value = 7866
$ws.Range("M31").Value = -1174.2307
$ws.Range("N31").Value = -8456
$ws.Range("H34").Value = 2974.353
$ws.Range("I34").Value = 1469.2307
$ws.Range("J34").Value = 7866
$ws.Range("K34").Value = 1469.2307
$ws.Range("L34").Value = 7866
$ws.Range("M34").Value = -1267.2307
$ws.Range("N34").Value = -8270
$ws.Range("H58").Value = 2589.6667
$ws.Range("I58").Value = 1627.4445
$ws.Range("J58").Value = 4033
$ws.Range("K58").Value = 1627.4445
$ws.Range("L58").Value = 4033
$ws.Range("M58").Value = -1424.4445
$ws.Range("N58").Value = -4439
$ws.Range("H86").Value = 8898.4
$ws.Range("I86").Value = 9664.333000000001
$ws.Range("J86").Value = 7749.5
$ws.Range("K86").Value = 9664.333000000001
$ws.Range("L86").Value = 7749.5
$ws.Range("M86").Value = -8541.333000000001
$ws.Range("N86").Value = -9995.5
$ws.Range("H89").Value = 8898.4
$ws.Range("I89").Value = 9664.333000000001
$ws.Range("J89").Value = 7749.5
$ws.Range("K89").Value = 48321.665
$ws.Range("L89").Value = 38747.5
$ws.Range("M89").Value = -42705.665
$ws.Range("N89").Value = -49979.5
$ws.Range("H122").Value = 3747.4211
$ws.Range("I122").Value = 4589.9165
$ws.Range("J122").Value = 2303.1428
$ws.Range("K122").Value = 13769.7495
$ws.Range("L122").Value = 6909.428400000001
$ws.Range("M122").Value = -11319.7495
$ws.Range("N122").Value = -11809.4284
$ws.Range("H136").Value = 2589.6667
$ws.Range("I136").Value = 1627.4445
$ws.Range("J136").Value = 4033
$ws.Range("K136").Value = 4882.333500000001
$ws.Range("L136").Value = 12099
$ws.Range("M136").Value = -2332.333500000001
$ws.Range("N136").Value = -17199

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5992.9165
$ws.Range("I80").Value = 6057.3335
$ws.Range("J80").Value = 5799.6665
$ws.Range("K80").Value = 18172.0005
$ws.Range("L80").Value = 17398.9995
$ws.Range("M80").Value = -17236.0005
$ws.Range("N80").Value = -19270.9995
$ws.Range("H83").Value = 5992.9165
$ws.Range("I83").Value = 6057.3335
$ws.Range("J83").Value = 5799.6665
$ws.Range("K83").Value = 54516.0015
$ws.Range("L83").Value = 52196.9985
$ws.Range("M83").Value = -49836.0015
$ws.Range("N83").Value = -61556.9985
$ws.Range("H131").Value = 1915.6666
$ws.Range("I131").Value = 1833.3334
$ws.Range("J131").Value = 1998
$ws.Range("K131").Value = 5500.0002
$ws.Range("L131").Value = 5994
$ws.Range("M131").Value = -460.0002000000004
$ws.Range("N131").Value = -16074

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 775
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 775
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 775
$ws.Range("N5").Value = -999
$ws.Range("H113").Value = 2299.8
$ws.Range("I113").Value = 1250
$ws.Range("J113").Value = 2999.6667
$ws.Range("K113").Value = 1250
$ws.Range("L113").Value = 2999.6667
$ws.Range("M113").Value = 920
$ws.Range("N113").Value = -7339.6667
$ws.Range("H122").Value = 2068.0417
$ws.Range("I122").Value = 1943.6471
$ws.Range("J122").Value = 2370.1428
$ws.Range("K122").Value = 5830.9413
$ws.Range("L122").Value = 7110.428400000001
$ws.Range("M122").Value = -3380.9413
$ws.Range("N122").Value = -12010.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3323.25
$ws.Range("I7").Value = 3465.6667
$ws.Range("J7").Value = 2896
$ws.Range("K7").Value = 3465.6667
$ws.Range("L7").Value = 2896
$ws.Range("M7").Value = -3353.6667
$ws.Range("N7").Value = -3120
$ws.Range("H22").Value = 819.3333
$ws.Range("I22").Value = 600.5454999999999
$ws.Range("J22").Value = 1060
$ws.Range("K22").Value = 600.5454999999999
$ws.Range("L22").Value = 1060
$ws.Range("M22").Value = -305.5454999999999
$ws.Range("N22").Value = -1650
$ws.Range("H27").Value = 819.3333
$ws.Range("I27").Value = 600.5454999999999
$ws.Range("J27").Value = 1060
$ws.Range("K27").Value = 600.5454999999999
$ws.Range("L27").Value = 1060
$ws.Range("M27").Value = -493.5454999999999
$ws.Range("N27").Value = -1274
$ws.Range("H40").Value = 1717
$ws.Range("I40").Value = 1948.75
$ws.Range("J40").Value = 790
$ws.Range("K40").Value = 1948.75
$ws.Range("L40").Value = 790
$ws.Range("M40").Value = -1812.75
$ws.Range("N40").Value = -1062
$ws.Range("H94").Value = 42332.332
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 42332.332
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 42332.332
$ws.Range("N94").Value = -43684.332
$ws.Range("H95").Value = 31224.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 31224.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 31224.5
$ws.Range("N95").Value = -36716.5
$ws.Range("H99").Value = 16749.666
$ws.Range("I99").Value = 18375
$ws.Range("J99").Value = 13499
$ws.Range("K99").Value = 18375
$ws.Range("L99").Value = 13499
$ws.Range("M99").Value = -15380
$ws.Range("N99").Value = -19489
$ws.Range("H126").Value = 3323.25
$ws.Range("I126").Value = 3465.6667
$ws.Range("J126").Value = 2896
$ws.Range("K126").Value = 10397.0001
$ws.Range("L126").Value = 8688
$ws.Range("M126").Value = -7927.000100000001
$ws.Range("N126").Value = -13628
$ws.Range("H132").Value = 1631.4
$ws.Range("I132").Value = 1679.7222
$ws.Range("J132").Value = 1196.5
$ws.Range("K132").Value = 5039.1666
$ws.Range("L132").Value = 3589.5
$ws.Range("M132").Value = -2509.1666
$ws.Range("N132").Value = -8649.5
$ws.Range("H136").Value = 25002594
$ws.Range("I136").Value = 2704.3333
$ws.Range("J136").Value = 62502428
$ws.Range("K136").Value = 8112.999899999999
$ws.Range("L136").Value = 187507284
$ws.Range("M136").Value = -5562.999899999999
$ws.Range("N136").Value = -187512384

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 28812.25
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 28812.25
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 28812.25
$ws.Range("N46").Value = -29274.25
$ws.Range("H96").Value = 2000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 2000
$ws.Range("N96").Value = -4746
$ws.Range("H101").Value = 27444.334
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 27444.334
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 27444.334
$ws.Range("N101").Value = -33934.334
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 2948.4783
$ws.Range("I122").Value = 3069.3157
$ws.Range("J122").Value = 2374.5
$ws.Range("K122").Value = 9207.947100000001
$ws.Range("L122").Value = 7123.5
$ws.Range("M122").Value = -6757.947100000001
$ws.Range("N122").Value = -12023.5
$ws.Range("H126").Value = 3337.5186
$ws.Range("I126").Value = 2879.7917
$ws.Range("J126").Value = 6999.3335
$ws.Range("K126").Value = 8639.375100000001
$ws.Range("L126").Value = 20998.0005
$ws.Range("M126").Value = -6169.375100000001
$ws.Range("N126").Value = -25938.0005
$ws.Range("H130").Value = 43999.332
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 43999.332
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 43999.332
$ws.Range("N130").Value = -54039.332
$ws.Range("H132").Value = 3130.4119
$ws.Range("I132").Value = 2749.3872
$ws.Range("J132").Value = 7067.6665
$ws.Range("K132").Value = 8248.161599999999
$ws.Range("L132").Value = 21202.9995
$ws.Range("M132").Value = -5718.161599999999
$ws.Range("N132").Value = -26262.9995
$ws.Range("H134").Value = 28812.25
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 28812.25
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 86436.75
$ws.Range("N134").Value = -91506.75
